$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Step 1: bump the header date shown in the merged A1:D1 banner
$ws.Range("A1").Value = 45309

# Step 2: update the two price cells in the pricing table
$ws.Range("D29").Value = 19600
$ws.Range("D30").Value = 21660
